# Convert a 6-hex-digit RRGGBB string into the COM "RGB" integer
# (PowerPoint/VBA packs colors as 0x00BBGGRR, i.e. R + G*256 + B*65536).
function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# EMU -> point conversion (1 pt = 12700 EMU) used throughout for Left/Top/Width/Height.
function Emu($v) { return $v / 12700.0 }

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Drop slides 3-6 (the presentation keeps only the first two slides).
# ---------------------------------------------------------------------------
while ($p.Slides.Count -gt 2) {
    $p.Slides.Item($p.Slides.Count).Delete()
}

# ---------------------------------------------------------------------------
# 2. Slide 1 ("Ocean Blue" title slide) restyle.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape 2 "Rectangle 2": was the thin top band -> becomes the thick bottom band.
$rect2 = $s1.Shapes.Item(2)
$rect2.Left = Emu(0)
$rect2.Top = Emu(5029200)
$rect2.Width = Emu(12191695)
$rect2.Height = Emu(1828800)

# New shape "Rectangle 3": thin accent band sitting just above the bottom band.
# Duplicate Rectangle 2 so it inherits the identical style/fill/line formatting,
# then move it into place and fix up its name + z-order.
$rect3 = $rect2.Duplicate()
$rect3.Name = "Rectangle 3"
$rect3.Left = Emu(0)
$rect3.Top = Emu(4846320)
$rect3.Width = Emu(12191695)
$rect3.Height = Emu(365760)
$rect3.Fill.ForeColor.RGB = HexToRgb("0096C8")
$rect3.ZOrder(1)   # msoSendToBack
$rect3.ZOrder(2)   # msoBringForward
$rect3.ZOrder(2)   # msoBringForward -> lands right after Rectangle 2

# Grab both text boxes by their *original* names before either gets renamed
# (renaming first would make a later Item("TextBox 4") lookup ambiguous).
$title = $s1.Shapes.Item("TextBox 3")
$subtitle = $s1.Shapes.Item("TextBox 4")

# Title textbox ("Ocean Blue"): widen, bump size, switch font.
$title.Width = Emu(9144000)
$title.TextFrame.TextRange.Font.Size = 48
$title.TextFrame.TextRange.Font.Name = "Calibri"
$title.Name = "TextBox 4"

# Subtitle textbox: widen, recolor, switch font, change copy.
$subtitle.Width = Emu(9144000)
$subtitle.TextFrame.TextRange.Text = "Calm & Professional"
$subtitle.TextFrame.TextRange.Font.Color.RGB = HexToRgb("506478")
$subtitle.TextFrame.TextRange.Font.Name = "Calibri"
$subtitle.Name = "TextBox 5"

# Old small underline accent rectangle is gone in the new design.
$s1.Shapes.Item("Rectangle 5").Delete()

# ---------------------------------------------------------------------------
# 3. Slide 2 ("Key Features" -> "Blue Features") restyle.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$bg = $s2.Shapes.Item("Rectangle 1")
$bg.Fill.ForeColor.RGB = HexToRgb("EBF5FF")

$header = $s2.Shapes.Item("TextBox 3")
$header.Top = Emu(228600)
$header.TextFrame.TextRange.Text = "Blue Features"
$header.TextFrame.TextRange.Font.Name = "Calibri"

$card = $s2.Shapes.Item("Rounded Rectangle 4")
$card.Fill.ForeColor.RGB = HexToRgb("FFFFFF")

$body = $s2.Shapes.Item("TextBox 5")
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "• Ocean blue palette`r• Calm professional look`r• Healthcare ready`r• Trustworthy design`r• Corporate friendly"
$bodyTr.Font.Size = 20
$bodyTr.Font.Name = "Calibri"
